$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "village" / "block" columns (header row + single data row), then
# the "land_area" column — typed in this order so the shared-string
# table is built in the same sequence as the source edit.
$ws.Range("C1").Value = "village"
$ws.Range("D1").Value = "block"
$ws.Range("C2").Value = "asdad"
$ws.Range("D2").Value = "asdsa"
$ws.Range("E1").Value = "land_area"
$ws.Range("E2").Value = 40

# Match the column widths the author left on the new columns.
$ws.Columns.Item(3).ColumnWidth = 11.333333333333332
$ws.Columns.Item(4).ColumnWidth = 13

# Leave the selection on the last cell that was edited.
$ws.Range("E2").Select() | Out-Null
